$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "'11"
$ws.Range("B2").Value = $false
$ws.Range("C2").Value = 'The abstract discusses curcumin''s biological activities and its potential health benefits but does not focus on NAFLD or extra-hepatic cancer outcomes. It lacks a defined population, intervention, comparison, or relevant outcomes, making it irrelevant to the PICOS criteria.'
$ws.Range("D2").Value = 'The abstract does not specify a population related to NAFLD or cancer.'
$ws.Range("E2").Value = 'There is no mention of an intervention related to NAFLD management or observation.'
$ws.Range("F2").Value = 'No comparison group involving NAFLD or general population is described.'
$ws.Range("G2").Value = 'Outcomes related to extra-hepatic cancers or NAFLD are not addressed.'
$ws.Range("H2").Value = 'The study design is a review, not a retrospective cohort study.'

# Row 3
$ws.Range("A3").Value = "'12"
$ws.Range("B3").Value = $false
$ws.Range("C3").Value = 'While the abstract mentions intestinal fungi and their association with diseases including NAFLD, it does not specifically address extra-hepatic cancer outcomes or provide details on interventions or comparisons. The study design is also a review, not a retrospective cohort study.'
$ws.Range("D3").Value = 'The population includes individuals with various diseases but does not focus specifically on NAFLD patients.'
$ws.Range("E3").Value = 'There is no specific intervention related to NAFLD management or observation.'
$ws.Range("F3").Value = 'No comparison group involving NAFLD or general population is described.'
$ws.Range("G3").Value = 'Outcomes related to extra-hepatic cancers are not addressed.'
$ws.Range("H3").Value = 'The study design is a review, not a retrospective cohort study.'

# Row 4
$ws.Range("A4").Value = "'13"
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 'The abstract focuses on LMNA variants and multisystem progeroid syndrome, which is unrelated to NAFLD or extra-hepatic cancer outcomes. The population, intervention, and outcomes do not align with the PICOS criteria.'
$ws.Range("D4").Value = 'The population involves patients with LMNA variants, not NAFLD.'
$ws.Range("E4").Value = 'There is no intervention related to NAFLD management or observation.'
$ws.Range("F4").Value = 'No comparison group involving NAFLD or general population is described.'
$ws.Range("G4").Value = 'Outcomes related to extra-hepatic cancers are not addressed.'
$ws.Range("H4").Value = 'The study design is a case series, not a retrospective cohort study.'

# Row 5
$ws.Range("A5").Value = "'14"
$ws.Range("B5").Value = $false
$ws.Range("C5").Value = 'The abstract discusses sex and gender disparities in disease but does not focus on NAFLD or extra-hepatic cancer outcomes. It lacks a defined population, intervention, comparison, or relevant outcomes, making it irrelevant to the PICOS criteria.'
$ws.Range("D5").Value = 'The population is not specified as NAFLD patients.'
$ws.Range("E5").Value = 'There is no intervention related to NAFLD management or observation.'
$ws.Range("F5").Value = 'No comparison group involving NAFLD or general population is described.'
$ws.Range("G5").Value = 'Outcomes related to extra-hepatic cancers are not addressed.'
$ws.Range("H5").Value = 'The study design is a review, not a retrospective cohort study.'

# Row 6
$ws.Range("A6").Value = "'15"
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 'The abstract compares the safety of two drugs in premenopausal breast cancer patients, which is unrelated to NAFLD or extra-hepatic cancer outcomes. The population, intervention, and outcomes do not align with the PICOS criteria.'
$ws.Range("D6").Value = 'The population involves premenopausal breast cancer patients, not NAFLD patients.'
$ws.Range("E6").Value = 'The intervention involves drug treatment for breast cancer, not NAFLD management.'
$ws.Range("F6").Value = 'No comparison group involving NAFLD or general population is described.'
$ws.Range("G6").Value = 'Outcomes related to extra-hepatic cancers are not addressed.'
$ws.Range("H6").Value = 'The study design is a randomized clinical trial, not a retrospective cohort study.'

# Row 7
$ws.Range("A7").Value = "'16"
$ws.Range("B7").Value = $false
$ws.Range("C7").Value = 'The abstract reviews liver diseases during pregnancy, including NAFLD, but does not focus on extra-hepatic cancer outcomes or provide details on interventions or comparisons. The study design is a review, not a retrospective cohort study.'
$ws.Range("D7").Value = 'The population includes pregnant women with liver diseases, not specifically NAFLD patients.'
$ws.Range("E7").Value = 'There is no specific intervention related to NAFLD management or observation.'
$ws.Range("F7").Value = 'No comparison group involving NAFLD or general population is described.'
$ws.Range("G7").Value = 'Outcomes related to extra-hepatic cancers are not addressed.'
$ws.Range("H7").Value = 'The study design is a review, not a retrospective cohort study.'

# Row 8
$ws.Range("A8").Value = "'17"
$ws.Range("B8").Value = $false
$ws.Range("C8").Value = 'The abstract discusses Faecalibacterium prausnitzii and its role in intestinal health but does not focus on NAFLD or extra-hepatic cancer outcomes. It lacks a defined population, intervention, comparison, or relevant outcomes, making it irrelevant to the PICOS criteria.'
$ws.Range("D8").Value = 'The population is not specified as NAFLD patients.'
$ws.Range("E8").Value = 'There is no intervention related to NAFLD management or observation.'
$ws.Range("F8").Value = 'No comparison group involving NAFLD or general population is described.'
$ws.Range("G8").Value = 'Outcomes related to extra-hepatic cancers are not addressed.'
$ws.Range("H8").Value = 'The study design is a review, not a retrospective cohort study.'

# Row 9
$ws.Range("A9").Value = "'18"
$ws.Range("B9").Value = $false
$ws.Range("C9").Value = 'The abstract focuses on cellular senescence in chronic kidney disease, which is unrelated to NAFLD or extra-hepatic cancer outcomes. The population, intervention, and outcomes do not align with the PICOS criteria.'
$ws.Range("D9").Value = 'The population involves patients with chronic kidney disease, not NAFLD patients.'
$ws.Range("E9").Value = 'There is no intervention related to NAFLD management or observation.'
$ws.Range("F9").Value = 'No comparison group involving NAFLD or general population is described.'
$ws.Range("G9").Value = 'Outcomes related to extra-hepatic cancers are not addressed.'
$ws.Range("H9").Value = 'The study design is a review, not a retrospective cohort study.'

# Row 10
$ws.Range("A10").Value = "'10"

# Row 11
$ws.Range("A11").Value = "'0"
$ws.Range("B11").Value = $false
$ws.Range("C11").Value = 'The abstract discusses moderate alcohol consumption''s effects on NAFLD but does not address the incidence of extra-hepatic cancers or compare NAFLD patients with a non-NAFLD population. The focus is on liver fibrosis progression rather than cancer outcomes.'
$ws.Range("D11").Value = 'Patients with NAFLD are mentioned, but the abstract focuses on alcohol consumption effects rather than cancer risks in NAFLD specifically.'
$ws.Range("E11").Value = 'The intervention pertains to moderate alcohol consumption and its impact on liver fibrosis, not observation or management of NAFLD related to cancer risks.'
$ws.Range("F11").Value = 'No comparison group involving non-NAFLD patients or general population is described in relation to cancer outcomes.'
$ws.Range("G11").Value = 'The outcome focuses on liver fibrosis progression rather than the incidence of extra-hepatic cancers.'
$ws.Range("H11").Value = 'The study design description is unclear; it appears to be a narrative review rather than a retrospective cohort study.'

# Row 12
$ws.Range("A12").Value = "'1"
$ws.Range("B12").Value = $false
$ws.Range("C12").Value = 'This abstract provides an overview of pediatric liver diseases, including NAFLD, but does not investigate extra-hepatic cancers or include a comparison group. It lacks relevance to the PICOS criteria for cancer incidence in NAFLD patients.'
$ws.Range("D12").Value = 'The population includes children with liver diseases, but there is no specific mention of NAFLD patients in the context of cancer risks.'
$ws.Range("E12").Value = 'The intervention involves diagnosing liver disease in children, which is unrelated to managing NAFLD or observing cancer risks.'
$ws.Range("F12").Value = 'No comparison group involving non-NAFLD patients or general population is described.'
$ws.Range("G12").Value = 'The outcomes focus on liver function tests and diagnostic methods, not extra-hepatic cancer incidence.'
$ws.Range("H12").Value = 'The study design appears to be a narrative review rather than a retrospective cohort study.'

# Row 13
$ws.Range("A13").Value = "'2"
$ws.Range("B13").Value = $false
$ws.Range("C13").Value = 'The abstract explores transporter alterations in gastrointestinal and kidney functions due to liver dysfunction but does not address NAFLD patients'' cancer risks or provide a relevant comparison group. It lacks alignment with the PICOS criteria.'
$ws.Range("D13").Value = 'The population involves patients with liver dysfunction, but NAFLD-specific patients are not highlighted in the context of cancer risks.'
$ws.Range("E13").Value = 'The intervention pertains to studying transporter alterations, which is unrelated to NAFLD management or observation of cancer risks.'
$ws.Range("F13").Value = 'No comparison group involving non-NAFLD patients or general population is described.'
$ws.Range("G13").Value = 'The outcomes focus on drug pharmacokinetics and transporter functions, not extra-hepatic cancer incidence.'
$ws.Range("H13").Value = 'The study design appears to be a literature review rather than a retrospective cohort study.'

# Row 14
$ws.Range("A14").Value = "'3"
$ws.Range("B14").Value = $true
$ws.Range("C14").Value = 'This abstract directly addresses extra-hepatic complications of NAFLD, including various cancers such as colorectal cancer. It mentions the need for collaborative care and screening methods, aligning well with the PICOS criteria despite lacking detailed study design information.'
$ws.Range("D14").Value = '-'
$ws.Range("E14").Value = '-'
$ws.Range("F14").Value = '-'
$ws.Range("G14").Value = '-'
$ws.Range("H14").Value = 'The study design is described as a narrative review, which may not fully meet the retrospective cohort requirement but still provides valuable insights into NAFLD-related cancer risks.'

# Row 15
$ws.Range("A15").Value = "'4"
$ws.Range("B15").Value = $false
$ws.Range("C15").Value = 'The abstract focuses on sphingolipids and ER stress mechanisms without addressing NAFLD patients'' cancer risks or providing a relevant comparison group. It lacks alignment with the PICOS criteria.'
$ws.Range("D15").Value = 'The population is not specified as NAFLD patients in the context of cancer risks.'
$ws.Range("E15").Value = 'The intervention involves studying sphingolipids and ER stress, which is unrelated to NAFLD management or observation of cancer risks.'
$ws.Range("F15").Value = 'No comparison group involving non-NAFLD patients or general population is described.'
$ws.Range("G15").Value = 'The outcomes focus on ER stress and sphingolipid roles, not extra-hepatic cancer incidence.'
$ws.Range("H15").Value = 'The study design appears to be a narrative review rather than a retrospective cohort study.'

# Row 16
$ws.Range("A16").Value = "'5"
$ws.Range("B16").Value = $false
$ws.Range("C16").Value = 'The abstract discusses epigenetic mechanisms and their role in gene expression regulation but does not address NAFLD patients'' cancer risks or provide a relevant comparison group. It lacks alignment with the PICOS criteria.'
$ws.Range("D16").Value = 'The population is not specified as NAFLD patients in the context of cancer risks.'
$ws.Range("E16").Value = 'The intervention involves studying epigenetic events, which is unrelated to NAFLD management or observation of cancer risks.'
$ws.Range("F16").Value = 'No comparison group involving non-NAFLD patients or general population is described.'
$ws.Range("G16").Value = 'The outcomes focus on epigenetic regulation, not extra-hepatic cancer incidence.'
$ws.Range("H16").Value = 'The study design appears to be a narrative review rather than a retrospective cohort study.'

# Row 17
$ws.Range("A17").Value = "'6"
$ws.Range("B17").Value = $true
$ws.Range("C17").Value = 'This abstract examines the association between NAFLD and colorectal polyps, including adenomas and cancers, using observational studies. It provides gender-specific risk analysis, aligning well with the PICOS criteria despite some heterogeneity in study types.'
$ws.Range("D17").Value = '-'
$ws.Range("E17").Value = '-'
$ws.Range("F17").Value = '-'
$ws.Range("G17").Value = '-'
$ws.Range("H17").Value = '-'

# Row 18
$ws.Range("A18").Value = "'9"
$ws.Range("B18").Value = $false
$ws.Range("C18").Value = 'The abstract discusses scutellarin''s pharmacological effects, including anti-tumor properties, but does not specifically address NAFLD patients'' cancer risks or provide a relevant comparison group. It lacks alignment with the PICOS criteria.'
$ws.Range("D18").Value = 'The population is not specified as NAFLD patients in the context of cancer risks.'
$ws.Range("E18").Value = 'The intervention involves studying scutellarin''s effects, which is unrelated to NAFLD management or observation of cancer risks.'
$ws.Range("F18").Value = 'No comparison group involving non-NAFLD patients or general population is described.'
$ws.Range("G18").Value = 'The outcomes focus on scutellarin''s mechanisms, not extra-hepatic cancer incidence.'
$ws.Range("H18").Value = 'The study design appears to be a narrative review rather than a retrospective cohort study.'

# Row 19
$ws.Range("A19").Value = "'7"
$ws.Range("B19").Value = $false
$ws.Range("C19").Value = 'Not processed - Empty abstract'
$ws.Range("D19").Value = 'not applicable'
$ws.Range("E19").Value = 'not applicable'
$ws.Range("F19").Value = 'not applicable'
$ws.Range("G19").Value = 'not applicable'
$ws.Range("H19").Value = 'not applicable'

# Row 20
$ws.Range("A20").Value = "'8"
